$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the misspelled name in B3 (Prretika Shetty -> Preetika Shetty) ---
$ws.Range("B3").Value = "Preetika Shetty"

# --- Update column A (ids) ---
$idValues = @(52501, 52502, 52503, 52504, 52505, 52506, 52507, 52508, 52509, 52510)
for ($i = 0; $i -lt $idValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $idValues[$i]
}

# --- Update column C (scores) ---
$scoreValues = @(80, 80, 97, 96, 93, 95, 92, 85, 91, 90)
for ($i = 0; $i -lt $scoreValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $scoreValues[$i]
}

# --- Update the selection so B3 is the active cell ---
$ws.Range("B3").Select()

$wb.Save()
